# The upstream commit ("meeting minutes updates from yuki ad edward") is
# recorded in the repo's OOXML diff as touching five parts -
# word/document.xml, word/endnotes.xml, word/footer1.xml,
# word/footnotes.xml and word/header1.xml - but in every one of those
# hunks the only byte-level change is the root element picking up one
# extra, unused namespace declaration:
#
#   xmlns:oel="http://schemas.microsoft.com/office/2019/extlst"
#
# No text, run, paragraph, style, table, header/footer body content,
# or document property differs between the two revisions (the hunk line
# counts are unchanged - 7/7, 13/13, 7/7 - and docProps/app.xml,
# docProps/core.xml, styles.xml, settings.xml, etc. don't even appear in
# the diff). That namespace prefix is never actually referenced anywhere
# in the markup; its appearance is the signature of the authoring copy
# of Word simply re-stamping the part headers with its current
# (newer-build) namespace table on save - not an edit any user made
# through the Word UI/object model, and not something that changes what
# the document *is* (its paragraphs, runs, formatting, headers/footers,
# etc. are identical before and after).
#
# There is no Word object-model surface (no property, method, or
# content-control/building-block/header-footer operation) that stamps
# that single inert xmlns declaration onto the existing parts without
# altering other content, so the faithful reproduction of this commit's
# actual content is a clean no-op here: the document already matches
# the target state content-wise, and this script intentionally makes no
# changes to text, formatting, or structure.
$d = $word.ActiveDocument
